# Apply the benchmark summary update: refresh several L/M column values
# (new parameter sets + new best run) and move the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated L/M values for rows 2-6 and 17-21 (only rows that have L/M data)
$updates = @(
    @{ Row = 2;  L = 0.63;  M = 0.672 },
    @{ Row = 3;  L = 0.626; M = 0.669 },
    @{ Row = 4;  L = 0.635; M = 0.675 },
    @{ Row = 5;  L = 0.722; M = 0.778 },
    @{ Row = 6;  L = 0.579; M = 0.607 },
    @{ Row = 17; L = 0.728; M = 0.791 },
    @{ Row = 18; L = 0.709; M = 0.782 },
    @{ Row = 19; L = 0.747; M = 0.801 },
    @{ Row = 20; L = 0.754; M = 0.827 },
    @{ Row = 21; L = 0.76;  M = 0.799 }
)

foreach ($u in $updates) {
    $ws.Range("L$($u.Row)").Value = $u.L
    $ws.Range("M$($u.Row)").Value = $u.M
}

# Move the active selection to N21 (as reflected in the sheetView selection)
$ws.Range("N21").Select()
